$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.950.12'
$ws.Range("E2").Value = '  -1.31%  '
$ws.Range("D3").Value = '1.638.14'
$ws.Range("E3").Value = '  -0.58%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("E6").Value = '  -0.35%  '
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("E8").Value = '  -0.74%  '
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("E10").Value = '  -1.95%  '
$ws.Range("E11").Value = '  +0.25%  '
$ws.Range("D12").Value = '1.865.16'
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("D14").Value = '1.648.89'
$ws.Range("E14").Value = '  -0.35%  '
$ws.Range("E15").Value = '  -1.36%  '
$ws.Range("D16").Value = '0.0₃0763'
$ws.Range("E16").Value = '  -0.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.99'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.93%  '
$ws.Range("D18").Value = '25.958.79'
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.81%  '
$ws.Range("E21").Value = '  -1.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.93'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.49%  '
$ws.Range("E23").Value = '  -0.91%  '
$ws.Range("E24").Value = '  +0.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.96'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.59%  '
$ws.Range("E26").Value = '  +0.36%  '
$ws.Range("E27").Value = '  +3.35%  '
$ws.Range("E28").Value = '  -1.83%  '
$ws.Range("E29").Value = '  -0.58%  '
$ws.Range("E30").Value = '  -0.85%  '
$ws.Range("E31").Value = '  -1.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.30'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.26%  '
$ws.Range("E34").Value = '  -4.57%  '
$ws.Range("E35").Value = '  +1.74%  '
$ws.Range("E36").Value = '  -1.58%  '
$ws.Range("D37").Value = '1.136.10'
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("E38").Value = '  -1.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.47'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.11%  '
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("E41").Value = '  +0.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.50'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.29'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.22%  '
$ws.Range("E44").Value = '  -0.56%  '
$ws.Range("D45").Value = '1.774.83'
$ws.Range("E45").Value = '  -0.59%  '
$ws.Range("E46").Value = '  +4.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.71'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("E48").Value = '  +2.84%  '
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("E50").Value = '  -1.08%  '
$ws.Range("E51").Value = '  -0.53%  '
